$wb = $excel.ActiveWorkbook

# Map of sheet index (1-based, matching workbook.xml sheet order) -> list of (row, newValue)
# Sheet 1 = "展览", Sheet 2 = "演出", Sheet 3 = "本地生活", Sheet 4 = "全部类型"

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 35  # F2: 34 -> 35
$ws.Cells.Item(5, 6).Value = 186  # F5: 185 -> 186
$ws.Cells.Item(6, 6).Value = 3802  # F6: 3800 -> 3802
$ws.Cells.Item(7, 6).Value = 187  # F7: 186 -> 187
$ws.Cells.Item(8, 6).Value = 114  # F8: 113 -> 114
$ws.Cells.Item(10, 6).Value = 81  # F10: 78 -> 81
$ws.Cells.Item(12, 6).Value = 673  # F12: 672 -> 673
$ws.Cells.Item(13, 6).Value = 167  # F13: 165 -> 167
$ws.Cells.Item(14, 6).Value = 932  # F14: 921 -> 932
$ws.Cells.Item(16, 6).Value = 231  # F16: 228 -> 231
$ws.Cells.Item(20, 6).Value = 87  # F20: 86 -> 87
$ws.Cells.Item(21, 6).Value = 3333  # F21: 3317 -> 3333
$ws.Cells.Item(22, 6).Value = 5671  # F22: 5658 -> 5671
$ws.Cells.Item(23, 6).Value = 36  # F23: 35 -> 36
$ws.Cells.Item(24, 6).Value = 20  # F24: 19 -> 20
$ws.Cells.Item(27, 6).Value = 40  # F27: 39 -> 40
$ws.Cells.Item(28, 6).Value = 3223  # F28: 3221 -> 3223
$ws.Cells.Item(30, 6).Value = 14  # F30: 12 -> 14
$ws.Cells.Item(31, 6).Value = 2426  # F31: 2421 -> 2426
$ws.Cells.Item(32, 6).Value = 567  # F32: 566 -> 567
$ws.Cells.Item(35, 6).Value = 190  # F35: 188 -> 190
$ws.Cells.Item(36, 6).Value = 252  # F36: 251 -> 252
$ws.Cells.Item(37, 6).Value = 344  # F37: 342 -> 344
$ws.Cells.Item(38, 6).Value = 108  # F38: 106 -> 108
$ws.Cells.Item(40, 6).Value = 884  # F40: 880 -> 884
$ws.Cells.Item(41, 6).Value = 10  # F41: 3 -> 10
$ws.Cells.Item(42, 6).Value = 44  # F42: 43 -> 44

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 90  # F2: 89 -> 90

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 35  # F2: 34 -> 35
$ws.Cells.Item(5, 6).Value = 186  # F5: 185 -> 186
$ws.Cells.Item(6, 6).Value = 3802  # F6: 3800 -> 3802
$ws.Cells.Item(7, 6).Value = 187  # F7: 186 -> 187
$ws.Cells.Item(8, 6).Value = 114  # F8: 113 -> 114
$ws.Cells.Item(10, 6).Value = 90  # F10: 89 -> 90
$ws.Cells.Item(11, 6).Value = 81  # F11: 78 -> 81
$ws.Cells.Item(13, 6).Value = 673  # F13: 672 -> 673
$ws.Cells.Item(14, 6).Value = 167  # F14: 165 -> 167
$ws.Cells.Item(15, 6).Value = 932  # F15: 921 -> 932
$ws.Cells.Item(17, 6).Value = 231  # F17: 228 -> 231
$ws.Cells.Item(21, 6).Value = 87  # F21: 86 -> 87
$ws.Cells.Item(22, 6).Value = 3333  # F22: 3317 -> 3333
$ws.Cells.Item(23, 6).Value = 5671  # F23: 5658 -> 5671
$ws.Cells.Item(24, 6).Value = 36  # F24: 35 -> 36
$ws.Cells.Item(25, 6).Value = 20  # F25: 19 -> 20
$ws.Cells.Item(28, 6).Value = 40  # F28: 39 -> 40
$ws.Cells.Item(29, 6).Value = 3223  # F29: 3221 -> 3223
$ws.Cells.Item(31, 6).Value = 14  # F31: 12 -> 14
$ws.Cells.Item(32, 6).Value = 2426  # F32: 2421 -> 2426
$ws.Cells.Item(33, 6).Value = 567  # F33: 566 -> 567
$ws.Cells.Item(36, 6).Value = 190  # F36: 188 -> 190
$ws.Cells.Item(37, 6).Value = 252  # F37: 251 -> 252
$ws.Cells.Item(38, 6).Value = 344  # F38: 342 -> 344
$ws.Cells.Item(39, 6).Value = 108  # F39: 106 -> 108
$ws.Cells.Item(41, 6).Value = 884  # F41: 880 -> 884
$ws.Cells.Item(42, 6).Value = 10  # F42: 3 -> 10
$ws.Cells.Item(43, 6).Value = 44  # F43: 43 -> 44
